$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dados")

# Remove the redundant "grandes regiões e unidades da federação" header row (row 6),
# which had no data of its own. Deleting the entire row shifts every row below it
# up by one (row 7 "norte" becomes row 6, ..., row 38 "distrito federal" becomes row 37).
$ws.Rows.Item(6).Delete()
